# Updated cryptos list on Wed Dec 27 18:26:36 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.126.40"
$ws.Range("E2").Value = "  +2.37%  "
$ws.Range("D3").Value = "2.347.53"
$ws.Range("E3").Value = "  +6.53%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "`'311.47"
$ws.Range("E5").Value = "  +5.75%  "
$ws.Range("D6").Value = "`'109.04"
$ws.Range("E6").Value = "  +2.70%  "
$ws.Range("E7").Value = "  +3.50%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "`'0.631"
$ws.Range("E9").Value = "  +7.16%  "
$ws.Range("D10").Value = "`'42.77"
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").Value = "`'0.0939"
$ws.Range("E11").Value = "  +3.96%  "
$ws.Range("D12").Value = "`'8.92"
$ws.Range("E12").Value = "  +3.29%  "
$ws.Range("E13").Value = "  +13.12%  "
$ws.Range("E14").Value = "  +2.46%  "
$ws.Range("D15").Value = "`'16.23"
$ws.Range("E15").Value = "  +9.69%  "
$ws.Range("D16").Value = "2.697.10"
$ws.Range("E16").Value = "  +6.26%  "
$ws.Range("D17").Value = "2.350.49"
$ws.Range("E17").Value = "  +6.04%  "
$ws.Range("D18").Value = "43.073.14"
$ws.Range("E18").Value = "  +2.27%  "
$ws.Range("E19").Value = "  +4.66%  "
$ws.Range("D20").Value = "`'7.29"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("D21").Value = "`'75.16"
$ws.Range("E21").Value = "  +4.37%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "`'2.56"
$ws.Range("E23").Value = "  +13.90%  "
$ws.Range("D24").Value = "`'252.05"
$ws.Range("E24").Value = "  +11.46%  "
$ws.Range("D25").Value = "`'9.05"
$ws.Range("E25").Value = "  +2.48%  "
$ws.Range("D26").Value = "`'12.03"
$ws.Range("E26").Value = "  +5.04%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "`'39.16"
$ws.Range("E28").Value = "  +4.31%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "`'2.27"
$ws.Range("E29").Value = "  +1.44%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "`'22.55"
$ws.Range("E30").Value = "  +9.12%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "`'173.61"
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("B32").Value = "WEMIXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").Value = "`'3.16"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "`'0.0923"
$ws.Range("E33").Value = "  +6.43%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "`'5.94"
$ws.Range("E34").Value = "  +9.23%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "`'0.133"
$ws.Range("E35").Value = "  +7.18%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "`'5.00"
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "`'0.0377"
$ws.Range("E37").Value = "  +5.33%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "`'4.10"
$ws.Range("E38").Value = "  -3.20%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "`'0.104"
$ws.Range("E39").Value = "  +3.18%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "`'2.70"
$ws.Range("E40").Value = "  +11.19%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").Value = "`'72.48"
$ws.Range("E41").Value = "  +3.91%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "`'1.49"
$ws.Range("E42").Value = "  +16.15%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "`'0.232"
$ws.Range("E43").Value = "  +2.80%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "`'1.00"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").Value = "`'12.75"
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D46").Value = "`'5.63"
$ws.Range("E46").Value = "  +4.95%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "`'9.35"
$ws.Range("E47").Value = "  +11.58%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "`'109.54"
$ws.Range("E48").Value = "  +7.67%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").Value = "`'1.30"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "`'0.101"
$ws.Range("E50").Value = "  +2.79%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "`'69.88"
$ws.Range("E51").Value = "  +5.66%  "
